$d = $word.ActiveDocument

# Paragraph 2 is the "!!! Server braucht noch mehr Funktionalität" paragraph
# that must be removed entirely (including its paragraph mark), and
# paragraph 3 is the following empty, underlined paragraph that should
# receive the new "In Overleaf ..." text.
$pDelete = $d.Paragraphs(2)
$pTarget = $d.Paragraphs(3)

$deleteRange = $d.Range($pDelete.Range.Start, $pTarget.Range.Start)
$deleteRange.Delete()

# After the delete, the formerly-empty underlined paragraph is now
# paragraph 2. Rebuild it (keeping its underline pPr) with the new runs,
# including the proofErr spell-check markers bracketing "Overleaf".
$target = $d.Paragraphs(2)
$ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$aumlaut = [char]0x00E4
$xml = "<w:p xmlns:w=""$ns""><w:pPr><w:rPr><w:u w:val=""single""/></w:rPr></w:pPr>" +
        "<w:r><w:t xml:space=""preserve"">In </w:t></w:r>" +
        "<w:proofErr w:type=""spellStart""/>" +
        "<w:r><w:t>Overleaf</w:t></w:r>" +
        "<w:proofErr w:type=""spellEnd""/>" +
        "<w:r><w:t xml:space=""preserve""> stehen alle drinnen mit Pseudowerten, allerdings sind keine Eintr" + $aumlaut + "ge davon schon im Glossar, das fehlt noch. </w:t></w:r>" +
        "</w:p>"
$target.Range.InsertXML($xml)
